$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header for column D
$ws.Cells.Item(2, 4).Value = "canonical SMILES"

# Column D repeats the "canonical isomeric SMILES" values from column C
for ($r = 3; $r -le 8; $r++) {
    $ws.Cells.Item($r, 4).Value = $ws.Cells.Item($r, 3).Value2
}

# Match the column width recorded in the target workbook
$ws.Columns.Item(4).ColumnWidth = 36
